# Updates cryptos list values/percentages per the commit diff.
# Numeric-looking "Price" strings (e.g. "309.38") are written through a
# scratch cell with a leading apostrophe + Copy/PasteSpecial(values) so
# Excel keeps them as literal text (matching the source inlineStr cells)
# instead of silently re-typing them as numbers (which would also mangle
# values like "0.0300" -> 0.03 or introduce float rounding noise).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$scratch = $ws.Range("Z1")

# Row 2
$ws.Range("D2").Value = '42.718.73'
$ws.Range("E2").Value = '  -0.98%  '
# Row 3
$ws.Range("D3").Value = '2.527.93'
$ws.Range("E3").Value = '  -2.37%  '
# Row 4
$scratch.Value = "'0.999"
$scratch.Copy() | Out-Null
$ws.Range("D4").PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null
$ws.Range("E4").Value = '  -0.03%  '
# Row 5
$scratch.Value = "'309.38"
$scratch.Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null
$ws.Range("E5").Value = '  -1.90%  '
# Row 6
$scratch.Value = "'101.57"
$scratch.Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null
$ws.Range("E6").Value = '  +4.09%  '
# Row 7
$scratch.Value = "'0.568"
$scratch.Copy() | Out-Null
$ws.Range("D7").PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null
$ws.Range("E7").Value = '  -1.42%  '
# Row 8
$ws.Range("E8").Value = '  +0.05%  '
# Row 9
$scratch.Value = "'0.527"
$scratch.Copy() | Out-Null
$ws.Range("D9").PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null
$ws.Range("E9").Value = '  -1.60%  '
# Row 10
$scratch.Value = "'35.97"
$scratch.Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null
$ws.Range("E10").Value = '  +1.06%  '
# Row 11
$scratch.Value = "'0.0805"
$scratch.Copy() | Out-Null
$ws.Range("D11").PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null
$ws.Range("E11").Value = '  -1.04%  '
# Row 12
$scratch.Value = "'7.32"
$scratch.Copy() | Out-Null
$ws.Range("D12").PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null
$ws.Range("E12").Value = '  -2.27%  '
# Row 13
$ws.Range("E13").Value = '  +0.34%  '
# Row 14
$ws.Range("D14").Value = '2.912.50'
$ws.Range("E14").Value = '  -2.51%  '
# Row 15
$scratch.Value = "'15.72"
$scratch.Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null
$ws.Range("E15").Value = '  +2.93%  '
# Row 16
$ws.Range("D16").Value = '2.586.59'
$ws.Range("E16").Value = '  +2.77%  '
# Row 17
$scratch.Value = "'0.808"
$scratch.Copy() | Out-Null
$ws.Range("D17").PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null
$ws.Range("E17").Value = '  -4.38%  '
# Row 18
$ws.Range("D18").Value = '42.692.09'
$ws.Range("E18").Value = '  -1.22%  '
# Row 19
$scratch.Value = "'6.77"
$scratch.Copy() | Out-Null
$ws.Range("D19").PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null
$ws.Range("E19").Value = '  -0.96%  '
# Row 20
$ws.Range("D20").Value = '0.0₃0952'
$ws.Range("E20").Value = '  -1.47%  '
# Row 21
$scratch.Value = "'12.10"
$scratch.Copy() | Out-Null
$ws.Range("D21").PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null
$ws.Range("E21").Value = '  -4.33%  '
# Row 22
$scratch.Value = "'69.33"
$scratch.Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null
$ws.Range("E22").Value = '  -0.23%  '
# Row 23
$scratch.Value = "'244.51"
$scratch.Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null
$ws.Range("E23").Value = '  -2.72%  '
# Row 24
$ws.Range("E24").Value = '  -2.26%  '
# Row 25
$ws.Range("E25").Value = '  -1.81%  '
# Row 26
$ws.Range("E26").Value = '  +0.02%  '
# Row 27
$scratch.Value = "'26.19"
$scratch.Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null
$ws.Range("E27").Value = '  -4.18%  '
# Row 28
$ws.Range("E28").Value = '  -3.73%  '
# Row 29
$scratch.Value = "'39.11"
$scratch.Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null
$ws.Range("E29").Value = '  -4.22%  '
# Row 30
$scratch.Value = "'10.15"
$scratch.Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null
$ws.Range("E30").Value = '  -1.34%  '
# Row 31
$scratch.Value = "'157.58"
$scratch.Copy() | Out-Null
$ws.Range("D31").PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null
$ws.Range("E31").Value = '  +0.49%  '
# Row 32
$scratch.Value = "'5.80"
$scratch.Copy() | Out-Null
$ws.Range("D32").PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null
$ws.Range("E32").Value = '  -0.96%  '
# Row 33
$scratch.Value = "'2.80"
$scratch.Copy() | Out-Null
$ws.Range("D33").PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null
$ws.Range("E33").Value = '  +11.53%  '
# Row 34
$scratch.Value = "'0.0788"
$scratch.Copy() | Out-Null
$ws.Range("D34").PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null
$ws.Range("E34").Value = '  -2.24%  '
# Row 35
$scratch.Value = "'2.62"
$scratch.Copy() | Out-Null
$ws.Range("D35").PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null
$ws.Range("E35").Value = '  -2.76%  '
# Row 36
$scratch.Value = "'2.04"
$scratch.Copy() | Out-Null
$ws.Range("D36").PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null
$ws.Range("E36").Value = '  -5.48%  '
# Row 37
$scratch.Value = "'3.20"
$scratch.Copy() | Out-Null
$ws.Range("D37").PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null
$ws.Range("E37").Value = '  -6.72%  '
# Row 38
$scratch.Value = "'18.17"
$scratch.Copy() | Out-Null
$ws.Range("D38").PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null
$ws.Range("E38").Value = '  -3.44%  '
# Row 39
$scratch.Value = "'0.112"
$scratch.Copy() | Out-Null
$ws.Range("D39").PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null
$ws.Range("E39").Value = '  -0.47%  '
# Row 40
$ws.Range("E40").Value = '  +0.17%  '
# Row 41
$ws.Range("E41").Value = '  +5.89%  '
# Row 42
$scratch.Value = "'21.98"
$scratch.Copy() | Out-Null
$ws.Range("D42").PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null
$ws.Range("E42").Value = '  -5.09%  '
# Row 43
$ws.Range("E43").Value = '  +0.01%  '
# Row 44
$ws.Range("B44").Value = 'NEARProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$scratch.Value = "'3.30"
$scratch.Copy() | Out-Null
$ws.Range("D44").PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null
$ws.Range("E44").Value = '  +2.00%  '
# Row 45
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$scratch.Value = "'0.0300"
$scratch.Copy() | Out-Null
$ws.Range("D45").PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null
$ws.Range("E45").Value = '  -1.47%  '
# Row 46
$ws.Range("D46").Value = '1.995.48'
$ws.Range("E46").Value = '  -0.81%  '
# Row 47
$scratch.Value = "'8.88"
$scratch.Copy() | Out-Null
$ws.Range("D47").PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null
$ws.Range("E47").Value = '  -1.25%  '
# Row 48
$ws.Range("D48").Value = '2.768.03'
$ws.Range("E48").Value = '  -2.53%  '
# Row 49
$scratch.Value = "'80.36"
$scratch.Copy() | Out-Null
$ws.Range("D49").PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null
$ws.Range("E49").Value = '  -3.50%  '
# Row 50
$ws.Range("E50").Value = '  -2.84%  '
# Row 51
$scratch.Value = "'72.22"
$scratch.Copy() | Out-Null
$ws.Range("D51").PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null
$ws.Range("E51").Value = '  -2.85%  '
